$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.309.87'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '3.308.10'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '190.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '561.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -2.44%  '
$ws.Range("D9").Value = '3.299.10'
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.186'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.80%  '
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").Value = '3.837.69'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '614.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").Value = '66.293.44'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").Value = '3.311.00'
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.48%  '
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.21%  '
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.73'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '566.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("D37").Value = '3.735.56'
$ws.Range("E37").Value = '  -4.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.13%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '0.0₃0733'
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '34.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.34%  '
$ws.Range("E42").Value = '  -5.04%  '
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.82%  '
$ws.Range("E46").Value = '  -3.61%  '
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("E49").Value = '  -1.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.01%  '
